$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# The "Nathan 1999" etc. benefit column (E) held text placeholders
# ("–0,16 ", "–0,11 ", ...) that pointed at shared-string entries; the
# updated catalog instead carries the real negative numeric benefit
# values computed for those rows. Replace each one with its number so
# the cells become numeric (t="s" is dropped) and the now-unused shared
# strings disappear from xl/sharedStrings.xml.
$ws.Range("E8").Value  = -0.16
$ws.Range("E14").Value = -0.11
$ws.Range("E23").Value = -0.19
$ws.Range("E27").Value = -0.22
$ws.Range("E35").Value = -0.011
$ws.Range("E42").Value = -0.09
$ws.Range("E49").Value = -0.1
$ws.Range("E63").Value = -0.21
$ws.Range("E65").Value = -0.1

# Re-point the saved view: scroll the window down so row 40 is at the
# top and move the active selection to H72.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H72").Select()
